$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W8)
$ws1.Range("D2").Value = 87
$ws1.Range("H2").Value = 4.91
$ws1.Range("L2").Value = 1.15

# Row 3 (W9)
$ws1.Range("D3").Value = 90
$ws1.Range("H3").Value = 3.78
$ws1.Range("L3").Value = 1.18

# Row 4 (W10)
$ws1.Range("D4").Value = 89
$ws1.Range("H4").Value = 2.79
$ws1.Range("L4").Value = 1

# Row 5 (W11)
$ws1.Range("D5").Value = 85
$ws1.Range("H5").Value = 1.87
$ws1.Range("L5").Value = 1.06

# Row 6 (W12)
$ws1.Range("D6").Value = 88
$ws1.Range("H6").Value = 0.84
$ws1.Range("J6").Value = "Urgent"
$ws1.Range("L6").Value = 0.91

# Row 7 (W13)
$ws1.Range("D7").Value = 83
$ws1.Range("H7").Value = 0
$ws1.Range("L7").Value = 1.13

# Row 8 (W14)
$ws1.Range("L8").Value = 0.99

# Row 9 (W15)
$ws1.Range("L9").Value = 1.15

# Row 10 (W16)
$ws1.Range("L10").Value = 1.03

# Row 11 (W17)
$ws1.Range("L11").Value = 1.16

# Row 12 (W18)
$ws1.Range("L12").Value = 1.2

# Row 13 (W19)
$ws1.Range("L13").Value = 0.88

# Row 14 (W20)
$ws1.Range("L14").Value = 0.8100000000000001

# Row 15 (W21)
$ws1.Range("L15").Value = 1.01

# Row 16 (W22)
$ws1.Range("L16").Value = 0.96

# Row 17 (W23)
$ws1.Range("L17").Value = 0.83

# --- Sheet 2: "Summary" ---
# These "Value" column cells hold numeric-looking text (stored as text in
# the source data, not numbers), so force a Text number format before
# assigning, otherwise Excel would auto-convert them to numeric values.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "1360"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "700"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "352"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "90"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "78"
